# edit.ps1
#
# Applies the commit's change to Demo.docx:
#   - Keeps paragraph 1 ("Welcome To Demo..!") unchanged.
#   - Inserts a brand-new paragraph right after it containing:
#       "Hello " + <spell-check proofErr wrap>Ayush</...> + " \u2026!"
#   - The "_GoBack" bookmark (last-edit-position marker) that used to sit
#     at the end of paragraph 1 now ends up at the end of the new,
#     second paragraph (matching real Word's behaviour of re-anchoring
#     _GoBack to the most recently edited spot).

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Figure out where paragraph 1's text ends (just before its
#    paragraph mark) -- this is where the new paragraph break goes.
# ------------------------------------------------------------------
$para1 = $d.Paragraphs(1)
$splitPos = $para1.Range.End - 1

# ------------------------------------------------------------------
# 2. Remove the existing "_GoBack" bookmark up front. We will recreate
#    it in the correct spot once the new paragraph exists.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3. Build the WordprocessingML for the new paragraph and insert it
#    via Range.InsertXML so that we get full control over the markup,
#    including the <w:proofErr/> spell-check markers around "Ayush".
#    A trailing sentinel character (a Private-Use-Area code point that
#    will never legitimately appear in the document) is appended
#    temporarily; it gives us a safe, unambiguous anchor point to
#    (re)plant the _GoBack bookmark after the fact, then we delete the
#    sentinel so the final text is exactly right.
# ------------------------------------------------------------------
$sentinel = "&#xE000;"

$newParaXml = '<?xml version="1.0" standalone="yes"?>' +
    '<?mso-application progid="Word.Document"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
    '<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
    '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
    '</Relationships></pkg:xmlData></pkg:part>' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p>' +
    '<w:r><w:t xml:space="preserve">Hello </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Ayush</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> &#8230;!' + $sentinel + '</w:t></w:r>' +
    '</w:p>' +
    '<w:sectPr/>' +
    '</w:body>' +
    '</w:document></pkg:xmlData></pkg:part>' +
    '</pkg:package>'

$insertAt = $d.Range($splitPos, $splitPos)
$insertAt.InsertXML($newParaXml)

# ------------------------------------------------------------------
# 4. Re-create the "_GoBack" bookmark at the end of the new paragraph
#    (right before its paragraph mark, i.e. right before the sentinel
#    char we just inserted), then strip the sentinel back out.
# ------------------------------------------------------------------
$bookmarkPos = $d.Content.End - 2
$d.Bookmarks.Add("_GoBack", $d.Range($bookmarkPos, $bookmarkPos)) | Out-Null

$sentinelRange = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$sentinelRange.Delete()
